# Update the marksheet's "Total" row with corrected/total marks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Total row (row 12) "Right" count and the Correct/Total summary text.
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 125
$ws.Range("E12").Value = "125/140"
